# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# Changes applied:
#  1. About!C1        : last-updated date bumped 2024-03-15 -> 2024-03-28 (serial 45366 -> 45379)
#  2. RAF-capacity!B24: hydrogen combustion turbine capacity-credit multiplier 0.3 -> 1
#  3. RAF-capacity!B25: hydrogen combined cycle capacity-credit multiplier 0.3 -> 1
#  4. RAF-capacity column A widened (new custom width, ~29 chars)
#  5. RAF-capacity becomes the active/selected sheet (was RAF-generation),
#     with cell B25 selected and the view zoomed to 80%.

$wb = $excel.ActiveWorkbook

# --- 1. Bump the "last updated" date stamp on the About sheet -------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- 2 & 3. Raise the hydrogen plant RAF-capacity multipliers to 1 --------
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- 4. Widen column A on RAF-capacity -------------------------------------
$wsCapacity.Columns.Item(1).ColumnWidth = 28.1

# --- 5. Make RAF-capacity the active sheet, select B25, zoom to 80% -------
$wsCapacity.Activate() | Out-Null
$wsCapacity.Range("B25").Select() | Out-Null
$excel.ActiveWindow.Zoom = 80
